# Applies the "add note to database internals" edit.
#
# The document is a long list of mostly-empty spacer paragraphs
# interspersed with real content. We edit from the bottom of the
# document upward so earlier paragraph indices stay stable while we
# work.

$d = $word.ActiveDocument

# --- 4) Insert the new "Note:" paragraph after paragraph 182 -----------
# (between the block of empty spacer paragraphs and "Article by
# Microsoft:" near the end of the document)
$p182 = $d.Paragraphs.Item(182)
$noteRange = $p182.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item(183)
$newPara.Format.LeftIndent = 0
$newPara.Range.Text = "Note: The database doesn" + [char]8217 + "t read a single row, it reads a page or more in a single IO and we get a lot of rows in that IO."
$newPara.Range.Font.Size = 14
$newPara.Range.Font.SizeBi = 14
$newPara.Range.Font.Underline = 0

# --- 3) Insert 3 blank paragraphs before paragraph 119 ------------------
# (extends the run of blank "ind left=0" spacer paragraphs just before
# the "Q: If we create a clustered index..." question)
$p118 = $d.Paragraphs.Item(118)
for ($i = 0; $i -lt 3; $i++) {
    $newRange = $p118.Range.InsertParagraphAfter()
}
for ($i = 119; $i -le 121; $i++) {
    $d.Paragraphs.Item($i).Format.LeftIndent = 0
}

# --- 2) Collapse 3 blank "ind left=720" spacer paragraphs into 1 --------
# (paragraphs 94-96, just before "Q: How The heap is stored on the
# disk?") blank "ind left=0" paragraph.
$d.Paragraphs.Item(95).Range.Delete()
$d.Paragraphs.Item(95).Range.Delete()
$d.Paragraphs.Item(94).Format.LeftIndent = 0

# --- 1) Delete 2 blank "ind left=0" spacer paragraphs -------------------
# (paragraphs 29-30, just before "Q: What is the difference between the
# clustered Index and the non-clustered Index?")
$d.Paragraphs.Item(29).Range.Delete()
$d.Paragraphs.Item(29).Range.Delete()
